$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 20126.8
$ws.Range("I9").Value = 20126.8
$ws.Range("K9").Value = 20126.8
$ws.Range("M9").Value = -19957.8
$ws.Range("H33").Value = 246.64706
$ws.Range("I33").Value = 145.41667
$ws.Range("K33").Value = 145.41667
$ws.Range("M33").Value = 83.58332999999999
$ws.Range("H51").Value = 4574.5
$ws.Range("H70").Value = 5958.3335
$ws.Range("J70").Value = 6650
$ws.Range("L70").Value = 19950
$ws.Range("N70").Value = -20490
$ws.Range("H73").Value = 5958.3335
$ws.Range("J73").Value = 6650
$ws.Range("L73").Value = 19950
$ws.Range("N73").Value = -21822
$ws.Range("H135").Value = 31250920
$ws.Range("I135").Value = 38462030
$ws.Range("K135").Value = 346158270
$ws.Range("M135").Value = -346155735
$ws.Range("H137").Value = 1923.4166
$ws.Range("I137").Value = 923.25
$ws.Range("K137").Value = 2769.75
$ws.Range("M137").Value = -219.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7179.1113
$ws.Range("I45").Value = 8171.143
$ws.Range("K45").Value = 8171.143
$ws.Range("M45").Value = -7794.143
$ws.Range("H61").Value = 71431400
$ws.Range("I61").Value = 90910230
$ws.Range("J61").Value = 8997
$ws.Range("K61").Value = 90910230
$ws.Range("L61").Value = 8997
$ws.Range("M61").Value = -90910018
$ws.Range("N61").Value = -9421
$ws.Range("H88").Value = 3374.5
$ws.Range("J88").Value = 3249.8
$ws.Range("L88").Value = 3249.8
$ws.Range("N88").Value = -4061.8
$ws.Range("H91").Value = 3374.5
$ws.Range("J91").Value = 3249.8
$ws.Range("L91").Value = 3249.8
$ws.Range("N91").Value = -6057.8
$ws.Range("H102").Value = 6698
$ws.Range("I102").Value = 7519.8887
$ws.Range("K102").Value = 7519.8887
$ws.Range("M102").Value = -5897.8887
$ws.Range("H136").Value = 71431400
$ws.Range("I136").Value = 90910230
$ws.Range("J136").Value = 8997
$ws.Range("K136").Value = 272730690
$ws.Range("L136").Value = 26991
$ws.Range("M136").Value = -272728140
$ws.Range("N136").Value = -32091
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4400.8
$ws.Range("I86").Value = 4400.8
$ws.Range("K86").Value = 4400.8
$ws.Range("M86").Value = -3277.8
$ws.Range("H89").Value = 4400.8
$ws.Range("I89").Value = 4400.8
$ws.Range("K89").Value = 22004
$ws.Range("M89").Value = -16388
$ws.Range("H94").Value = 3196.1
$ws.Range("I94").Value = 3196.1
$ws.Range("K94").Value = 3196.1
$ws.Range("M94").Value = -2745.1
$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524
$ws.Range("H134").Value = 31251432
$ws.Range("I134").Value = 31251432
$ws.Range("K134").Value = 93754296
$ws.Range("M134").Value = -93751761

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 26114.5
$ws.Range("I16").Value = 33819.332
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 33819.332
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -33532.332
$ws.Range("N16").Value = -3574
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H58").Value = 22735892
$ws.Range("I58").Value = 27787346
$ws.Range("K58").Value = 27787346
$ws.Range("M58").Value = -27787143
$ws.Range("H99").Value = 1922.5385
$ws.Range("I99").Value = 2019.5
$ws.Range("J99").Value = 1599.3334
$ws.Range("K99").Value = 2019.5
$ws.Range("L99").Value = 1599.3334
$ws.Range("M99").Value = -521.5
$ws.Range("N99").Value = -4595.3334
$ws.Range("H107").Value = 101132.9
$ws.Range("I107").Value = 389
$ws.Range("K107").Value = 389
$ws.Range("M107").Value = 1531
$ws.Range("H113").Value = 26114.5
$ws.Range("I113").Value = 33819.332
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 33819.332
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -31649.332
$ws.Range("N113").Value = -7340
$ws.Range("H122").Value = 1558.8
$ws.Range("I122").Value = 1561.5
$ws.Range("K122").Value = 4684.5
$ws.Range("M122").Value = -2234.5
$ws.Range("H126").Value = 1922.5385
$ws.Range("I126").Value = 2019.5
$ws.Range("J126").Value = 1599.3334
$ws.Range("K126").Value = 6058.5
$ws.Range("L126").Value = 4798.0002
$ws.Range("M126").Value = -3588.5
$ws.Range("N126").Value = -9738.0002
$ws.Range("H136").Value = 22735892
$ws.Range("I136").Value = 27787346
$ws.Range("K136").Value = 83362038
$ws.Range("M136").Value = -83359488

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 23400
$ws.Range("I106").Value = 2000
$ws.Range("J106").Value = 28750
$ws.Range("K106").Value = 6000
$ws.Range("L106").Value = 86250
$ws.Range("M106").Value = -5054
$ws.Range("N106").Value = -88142
$ws.Range("H117").Value = 2439.0908
$ws.Range("I117").Value = 1990
$ws.Range("J117").Value = 2484
$ws.Range("K117").Value = 5970
$ws.Range("L117").Value = 7452
$ws.Range("M117").Value = -2528
$ws.Range("N117").Value = -14336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2051.36
$ws.Range("I102").Value = 2072.1365
$ws.Range("K102").Value = 2072.1365
$ws.Range("M102").Value = -450.1365000000001
$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2086
$ws.Range("I40").Value = 1874.8572
$ws.Range("J40").Value = 2825
$ws.Range("K40").Value = 1874.8572
$ws.Range("L40").Value = 2825
$ws.Range("M40").Value = -1738.8572
$ws.Range("N40").Value = -3097
$ws.Range("H74").Value = 71515.2
$ws.Range("I74").Value = 61789.5
$ws.Range("K74").Value = 61789.5
$ws.Range("M74").Value = -60791.5
$ws.Range("H77").Value = 71515.2
$ws.Range("I77").Value = 61789.5
$ws.Range("K77").Value = 185368.5
$ws.Range("M77").Value = -180376.5
$ws.Range("H136").Value = 5049.357
$ws.Range("I136").Value = 4366.25
$ws.Range("K136").Value = 13098.75
$ws.Range("M136").Value = -10548.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 24851.5
$ws.Range("J104").Value = 24851.5
$ws.Range("L104").Value = 24851.5
$ws.Range("N104").Value = -31839.5
$ws.Range("H107").Value = 1034.8572
$ws.Range("I107").Value = 863
$ws.Range("J107").Value = 1264
$ws.Range("K107").Value = 2589
$ws.Range("L107").Value = 3792
$ws.Range("M107").Value = -669
$ws.Range("N107").Value = -7632
$ws.Range("H122").Value = 1282.8422
$ws.Range("I122").Value = 1282.8422
$ws.Range("K122").Value = 3848.5266
$ws.Range("M122").Value = -1398.5266

Write-Host "Applied all updates"